$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.185.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.645.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.00%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.644.07'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.626'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.707'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.159'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000288'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.205.38'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.637.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.33%  '
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.87%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.126'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.16%  '
$ws.Range("B19").Value = 'Polygon'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.87%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.886.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '407.62'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.40'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.41'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -12.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.69%  '
$ws.Range("E34").Value = '  -6.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '64.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '42.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -9.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '594.06'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0881'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -9.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.395'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.997'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.135'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.99'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0434'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -11.95%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.134'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.23%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -10.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.680.12'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.13%  '
